# Natmi following Dr Hou advice
# Updates Ligand/Receptor-expressing cell counts (E,K: 1 -> 3) and recomputes
# the dependent expression / specificity metrics for rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.760574
$ws.Range("H2").Value = 5.281722
$ws.Range("I2").Value = 0.2878920521313718
$ws.Range("J2").Value = 0.2878920521313718
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7065936666666666
$ws.Range("N2").Value = 2.119781
$ws.Range("O2").Value = 0.005187843618793344
$ws.Range("P2").Value = 0.005187843618793344
$ws.Range("Q2").Value = 1.244010438098
$ws.Range("R2").Value = 11.196093942882
$ws.Range("S2").Value = 0.001493538945551058
$ws.Range("T2").Value = 0.001493538945551058

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.760574
$ws.Range("H3").Value = 5.281722
$ws.Range("I3").Value = 0.2878920521313718
$ws.Range("J3").Value = 0.2878920521313718
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 111.9320066666667
$ws.Range("N3").Value = 335.79602
$ws.Range("O3").Value = 0.8218100075305903
$ws.Range("P3").Value = 0.8218100075305903
$ws.Range("Q3").Value = 197.06458070516
$ws.Range("R3").Value = 1773.58122634644
$ws.Range("S3").Value = 0.2365925695300798
$ws.Range("T3").Value = 0.2365925695300798

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.760574
$ws.Range("H4").Value = 5.281722
$ws.Range("I4").Value = 0.2878920521313718
$ws.Range("J4").Value = 0.2878920521313718
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 23.563205
$ws.Range("N4").Value = 70.689615
$ws.Range("O4").Value = 0.1730021488506163
$ws.Range("P4").Value = 0.1730021488506163
$ws.Range("Q4").Value = 41.48476607967
$ws.Range("R4").Value = 373.36289471703
$ws.Range("S4").Value = 0.04980594365574097
$ws.Range("T4").Value = 0.04980594365574098

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.1366213333333333
$ws.Range("H5").Value = 0.409864
$ws.Range("I5").Value = 0.02234055258015711
$ws.Range("J5").Value = 0.02234055258015711
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7065936666666666
$ws.Range("N5").Value = 2.119781
$ws.Range("O5").Value = 0.005187843618793344
$ws.Range("P5").Value = 0.005187843618793344
$ws.Range("Q5").Value = 0.09653576886488888
$ws.Range("R5").Value = 0.8688219197839999
$ws.Range("S5").Value = 0.0001158992931432852
$ws.Range("T5").Value = 0.0001158992931432852

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.1366213333333333
$ws.Range("H6").Value = 0.409864
$ws.Range("I6").Value = 0.02234055258015711
$ws.Range("J6").Value = 0.02234055258015711
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 111.9320066666667
$ws.Range("N6").Value = 335.79602
$ws.Range("O6").Value = 0.8218100075305903
$ws.Range("P6").Value = 0.8218100075305903
$ws.Range("Q6").Value = 15.29229999347556
$ws.Range("R6").Value = 137.63069994128
$ws.Range("S6").Value = 0.01835968968413646
$ws.Range("T6").Value = 0.01835968968413646

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.1366213333333333
$ws.Range("H7").Value = 0.409864
$ws.Range("I7").Value = 0.02234055258015711
$ws.Range("J7").Value = 0.02234055258015711
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 23.563205
$ws.Range("N7").Value = 70.689615
$ws.Range("O7").Value = 0.1730021488506163
$ws.Range("P7").Value = 0.1730021488506163
$ws.Range("Q7").Value = 3.219236484706667
$ws.Range("R7").Value = 28.97312836236
$ws.Range("S7").Value = 0.003864963602877361
$ws.Range("T7").Value = 0.003864963602877361

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.218201
$ws.Range("H8").Value = 12.654603
$ws.Range("I8").Value = 0.6897673952884711
$ws.Range("J8").Value = 0.6897673952884711
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7065936666666666
$ws.Range("N8").Value = 2.119781
$ws.Range("O8").Value = 0.005187843618793344
$ws.Range("P8").Value = 0.005187843618793344
$ws.Range("Q8").Value = 2.980554111326999
$ws.Range("R8").Value = 26.824987001943
$ws.Range("S8").Value = 0.003578405380099001
$ws.Range("T8").Value = 0.003578405380099001

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.218201
$ws.Range("H9").Value = 12.654603
$ws.Range("I9").Value = 0.6897673952884711
$ws.Range("J9").Value = 0.6897673952884711
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 111.9320066666667
$ws.Range("N9").Value = 335.79602
$ws.Range("O9").Value = 0.8218100075305903
$ws.Range("P9").Value = 0.8218100075305903
$ws.Range("Q9").Value = 472.15170245334
$ws.Range("R9").Value = 4249.36532208006
$ws.Range("S9").Value = 0.566857748316374
$ws.Range("T9").Value = 0.566857748316374

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.218201
$ws.Range("H10").Value = 12.654603
$ws.Range("I10").Value = 0.6897673952884711
$ws.Range("J10").Value = 0.6897673952884711
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.563205
$ws.Range("N10").Value = 70.689615
$ws.Range("O10").Value = 0.1730021488506163
$ws.Range("P10").Value = 0.1730021488506163
$ws.Range("Q10").Value = 99.39433489420499
$ws.Range("R10").Value = 894.549014047845
$ws.Range("S10").Value = 0.119331241591998
$ws.Range("T10").Value = 0.119331241591998
